# Commit: "Fix a small mistyped"
# The only real content change is a typo fix in the heading:
#   "WieFrame-дизайн" -> "WireFrame-дизайн"
$d = $word.ActiveDocument

$d.Content.Find.Execute("WieFrame", $true, $false, $false, $false, $false, $true, 1, $false, "WireFrame", 2)
